$d = $word.ActiveDocument

# Locate the author paragraph ("Vijay Panthayi") so the new source-code
# block is inserted right after it, regardless of its exact index.
$authorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Panthayi*") {
        $authorIndex = $i
        break
    }
}

$authorPara = $d.Paragraphs($authorIndex)

# Insert a brand new paragraph right after the author paragraph.
$authorPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs($authorIndex + 1)
$newFull = $newPara.Range

# Style the paragraph as a Pandoc "Source Code" block and fill in the
# knitr warning text.
$newFull.Style = "SourceCode"
$newFull.Text = "Warning: package 'knitr' was built under R version 4.2.3"

# Apply the "Verbatim Char" character style to the text run only (not
# the paragraph mark), matching the target markup's <w:rStyle>.
$textOnly = $d.Range($newFull.Start, $newFull.End - 1)
$textOnly.Style = "VerbatimChar"
